$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.362.50'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '1.874.06'
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.34'
$ws.Range("E5").Value = '  +0.48%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4782'
$ws.Range("E7").Value = '  -1.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2825'
$ws.Range("E8").Value = '  -2.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06525'
$ws.Range("E9").Value = '  -1.17%  '
$ws.Range("D10").Value = '1.872.73'
$ws.Range("E10").Value = '  -0.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07458'
$ws.Range("E11").Value = '  +1.67%  '
$ws.Range("E12").Value = '  -2.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.099'
$ws.Range("E13").Value = '  -1.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.15'
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6548'
$ws.Range("E15").Value = '  -1.02%  '
$ws.Range("D16").Value = '30.337.92'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("E17").Value = '  -1.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9996'
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007602'
$ws.Range("E19").Value = '  -2.34%  '
$ws.Range("D20").Value = '2.124.00'
$ws.Range("E20").Value = '  -1.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.305'
$ws.Range("E21").Value = '  -2.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '218.85'
$ws.Range("E23").Value = '  +13.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.210'
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.339'
$ws.Range("E25").Value = '  -0.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.23'
$ws.Range("E26").Value = '  +1.70%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.46'
$ws.Range("E27").Value = '  +1.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.972'
$ws.Range("E28").Value = '  +2.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.459'
$ws.Range("E29").Value = '  +0.36%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.316'
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09332'
$ws.Range("E31").Value = '  +1.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.035'
$ws.Range("E32").Value = '  -0.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05077'
$ws.Range("E33").Value = '  -0.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.204'
$ws.Range("E34").Value = '  +5.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7504'
$ws.Range("E35").Value = '  +3.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.715'
$ws.Range("E36").Value = '  +0.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01827'
$ws.Range("E37").Value = '  +2.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.615'
$ws.Range("E38").Value = '  -1.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.064'
$ws.Range("E39").Value = '  -0.40%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9056'
$ws.Range("E40").Value = '  -1.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '106.94'
$ws.Range("E41").Value = '  +0.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.896'
$ws.Range("E42").Value = '  +0.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4274'
$ws.Range("E43").Value = '  -1.01%  '
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.370'
$ws.Range("E45").Value = '  -1.39%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.36'
$ws.Range("E46").Value = '  -0.76%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1283'
$ws.Range("E47").Value = '  -3.57%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.471'
$ws.Range("E48").Value = '  -7.81%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.928'
$ws.Range("E49").Value = '  -0.19%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.60'
$ws.Range("E50").Value = '  -1.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3887'
$ws.Range("E51").Value = '  +0.66%  '
